$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 4336880.15
$ws.Cells.Item(2, 4).Value = 221.0612707332225

$ws.Cells.Item(3, 1).Value = "FL"
$ws.Cells.Item(3, 2).Value = 3663289.79
$ws.Cells.Item(3, 3).Value = 20598139
$ws.Cells.Item(3, 4).Value = 177.8456679994246
$ws.Cells.Item(3, 6).Value = 3

$ws.Cells.Item(4, 1).Value = "CA"
$ws.Cells.Item(4, 2).Value = 6256537
$ws.Cells.Item(4, 3).Value = 39148760
$ws.Cells.Item(4, 4).Value = 159.8144360127882
$ws.Cells.Item(4, 6).Value = 1

$ws.Cells.Item(5, 2).Value = 1314770
$ws.Cells.Item(5, 4).Value = 148.0289286741662

$ws.Cells.Item(9, 2).Value = 479750
$ws.Cells.Item(9, 4).Value = 48.17982205953952

$ws.Cells.Item(10, 1).Value = "WA"
$ws.Cells.Item(10, 2).Value = 281800
$ws.Cells.Item(10, 3).Value = 7294336
$ws.Cells.Item(10, 4).Value = 38.63271447874077
$ws.Cells.Item(10, 6).Value = 12

$ws.Cells.Item(11, 1).Value = "WI"
$ws.Cells.Item(11, 2).Value = 193691
$ws.Cells.Item(11, 3).Value = 5778394
$ws.Cells.Item(11, 4).Value = 33.51986728492381
$ws.Cells.Item(11, 6).Value = 19

$ws.Cells.Item(12, 1).Value = "PA"
$ws.Cells.Item(12, 2).Value = 402050
$ws.Cells.Item(12, 3).Value = 12791181
$ws.Cells.Item(12, 4).Value = 31.43181227753716
$ws.Cells.Item(12, 6).Value = 6

$ws.Cells.Item(13, 1).Value = "IL"
$ws.Cells.Item(13, 2).Value = 310700
$ws.Cells.Item(13, 3).Value = 12821497
$ws.Cells.Item(13, 4).Value = 24.23273974950039
$ws.Cells.Item(13, 6).Value = 5

$ws.Cells.Item(14, 1).Value = "DE"
$ws.Cells.Item(14, 2).Value = 17025
$ws.Cells.Item(14, 3).Value = 949495
$ws.Cells.Item(14, 4).Value = 17.9305841526285
$ws.Cells.Item(14, 6).Value = 40

$ws.Cells.Item(15, 1).Value = "OH"
$ws.Cells.Item(15, 2).Value = 181505
$ws.Cells.Item(15, 3).Value = 11641879
$ws.Cells.Item(15, 4).Value = 15.59069631285465
$ws.Cells.Item(15, 6).Value = 7

$ws.Cells.Item(16, 1).Value = "RI"
$ws.Cells.Item(16, 2).Value = 16025
$ws.Cells.Item(16, 3).Value = 1056611
$ws.Cells.Item(16, 4).Value = 15.16641412970336
$ws.Cells.Item(16, 6).Value = 38

$ws.Cells.Item(17, 1).Value = "MA"
$ws.Cells.Item(17, 2).Value = 82100
$ws.Cells.Item(17, 3).Value = 6830193
$ws.Cells.Item(17, 4).Value = 12.02015814194416
$ws.Cells.Item(17, 6).Value = 14

$ws.Cells.Item(18, 1).Value = "HI"
$ws.Cells.Item(18, 2).Value = 17000
$ws.Cells.Item(18, 3).Value = 1422029
$ws.Cells.Item(18, 4).Value = 11.95474916474981
$ws.Cells.Item(18, 6).Value = 35

$ws.Cells.Item(19, 1).Value = "AZ"
$ws.Cells.Item(19, 2).Value = 80800
$ws.Cells.Item(19, 3).Value = 6946685
$ws.Cells.Item(19, 4).Value = 11.63144722986576
$ws.Cells.Item(19, 6).Value = 13

$ws.Cells.Item(20, 1).Value = "CT"
$ws.Cells.Item(20, 2).Value = 36150
$ws.Cells.Item(20, 3).Value = 3581504
$ws.Cells.Item(20, 4).Value = 10.09352495487929
$ws.Cells.Item(20, 6).Value = 26

$ws.Cells.Item(28, 1).Value = "ME"
$ws.Cells.Item(28, 2).Value = 5900
$ws.Cells.Item(28, 3).Value = 1332813
$ws.Cells.Item(28, 4).Value = 4.426727530418747
$ws.Cells.Item(28, 6).Value = 37

$ws.Cells.Item(29, 1).Value = "NH"
$ws.Cells.Item(29, 2).Value = 5225
$ws.Cells.Item(29, 3).Value = 1343622
$ws.Cells.Item(29, 4).Value = 3.888742518357097
$ws.Cells.Item(29, 6).Value = 36

$ws.Cells.Item(30, 1).Value = "UT"
$ws.Cells.Item(30, 2).Value = 8250
$ws.Cells.Item(30, 3).Value = 3045350
$ws.Cells.Item(30, 4).Value = 2.709048221058335
$ws.Cells.Item(30, 6).Value = 28

$ws.Cells.Item(31, 1).Value = "KS"
$ws.Cells.Item(31, 2).Value = 7500
$ws.Cells.Item(31, 3).Value = 2908776
$ws.Cells.Item(31, 4).Value = 2.57840411224515
$ws.Cells.Item(31, 6).Value = 31

$ws.Cells.Item(32, 1).Value = "MO"
$ws.Cells.Item(32, 2).Value = 15050
$ws.Cells.Item(32, 3).Value = 6090062
$ws.Cells.Item(32, 4).Value = 2.471239209058956
$ws.Cells.Item(32, 6).Value = 17

$ws.Cells.Item(33, 1).Value = "OR"
$ws.Cells.Item(33, 2).Value = 8700
$ws.Cells.Item(33, 3).Value = 4081943
$ws.Cells.Item(33, 4).Value = 2.13133794371945
$ws.Cells.Item(33, 6).Value = 24

$ws.Cells.Item(34, 1).Value = "AR"
$ws.Cells.Item(34, 2).Value = 5000
$ws.Cells.Item(34, 3).Value = 2990671
$ws.Cells.Item(34, 4).Value = 1.671865611429676
$ws.Cells.Item(34, 6).Value = 29

$ws.Cells.Item(35, 1).Value = "TX"
$ws.Cells.Item(35, 2).Value = 44920
$ws.Cells.Item(35, 3).Value = 27885195
$ws.Cells.Item(35, 4).Value = 1.610890653624621
$ws.Cells.Item(35, 6).Value = 2

$ws.Cells.Item(36, 1).Value = "AL"
$ws.Cells.Item(36, 2).Value = 4000
$ws.Cells.Item(36, 3).Value = 4864680
$ws.Cells.Item(36, 4).Value = 0.8222534678540007
$ws.Cells.Item(36, 6).Value = 22

$ws.Cells.Item(37, 1).Value = "IA"
$ws.Cells.Item(37, 2).Value = 2100
$ws.Cells.Item(37, 3).Value = 3132499
$ws.Cells.Item(37, 4).Value = 0.6703912754640944
$ws.Cells.Item(37, 6).Value = 27

$ws.Cells.Item(38, 1).Value = "CO"
$ws.Cells.Item(38, 2).Value = 2750
$ws.Cells.Item(38, 3).Value = 5531141
$ws.Cells.Item(38, 4).Value = 0.4971849388760837
$ws.Cells.Item(38, 6).Value = 20

$ws.Cells.Item(39, 1).Value = "AK"
$ws.Cells.Item(39, 2).Value = 250
$ws.Cells.Item(39, 3).Value = 738516
$ws.Cells.Item(39, 4).Value = 0.3385167010599635
$ws.Cells.Item(39, 6).Value = 41
